{"js": "// Lattice multiplication exercises: replace each table cell's five lines\n// (problem, factor-digits line, separator, and the two lattice digit rows)\n// with new values, cell by cell, preserving the existing run formatting\n// (font size 32) and the \"<w:br/>\"-separated layout.\n//\n// The grid is row-major; \"after\" values below come directly from the\n// canonical OOXML diff applied to this document.\nconst afterGrid = [\n  [\"22 x 54\", \"  5    4\", \"  ----\", \"2|    |\", \"2|    |\"],\n  [\"94 x 27\", \"  2    7\", \"  ----\", \"9|    |\", \"4|    |\"],\n  [\"60 x 66\", \"  6    6\", \"  ----\", \"6|    |\", \"0|    |\"],\n  [\"10 x 42\", \"  4    2\", \"  ----\", \"1|    |\", \"0|    |\"],\n  [\"33 x 46\", \"  4    6\", \"  ----\", \"3|    |\", \"3|    |\"],\n  [\"42 x 71\", \"  7    1\", \"  ----\", \"4|    |\", \"2|    |\"],\n  [\"87 x 91\", \"  9    1\", \"  ----\", \"8|    |\", \"7|    |\"],\n  [\"20 x 80\", \"  8    0\", \"  ----\", \"2|    |\", \"0|    |\"],\n  [\"73 x 17\", \"  1    7\", \"  ----\", \"7|    |\", \"3|    |\"],\n  [\"33 x 15\", \"  1    5\", \"  ----\", \"3|    |\", \"3|    |\"],\n  [\"23 x 64\", \"  6    4\", \"  ----\", \"2|    |\", \"3|    |\"],\n  [\"65 x 59\", \"  5    9\", \"  ----\", \"6|    |\", \"5|    |\"],\n  [\"64 x 32\", \"  3    2\", \"  ----\", \"6|    |\", \"4|    |\"],\n  [\"99 x 48\", \"  4    8\", \"  ----\", \"9|    |\", \"9|    |\"],\n  [\"32 x 30\", \"  3    0\", \"  ----\", \"3|    |\", \"2|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet cellIndex = 0;\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const lines = afterGrid[cellIndex];\n    cellIndex++;\n    if (!lines) continue;\n\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    const newText = lines.join(\"\\u000b\");\n    // Replace the text of the cell's (only) paragraph while keeping its\n    // run-level formatting (sz 32) intact.\n    paragraphs.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Lattice multiplication exercises: replace each table cell's five lines\n# (problem, factor-digits line, separator, and the two lattice digit rows)\n# with new values, cell by cell, preserving the existing run formatting\n# (font size 32) and the vertical-tab-separated (\"<w:br/>\") layout.\n#\n# The grid is row-major; \"after\" values below come directly from the\n# canonical OOXML diff applied to this document.\n$afterGrid = @(\n    ,@('22 x 54', '  5    4', '  ----', '2|    |', '2|    |')\n    ,@('94 x 27', '  2    7', '  ----', '9|    |', '4|    |')\n    ,@('60 x 66', '  6    6', '  ----', '6|    |', '0|    |')\n    ,@('10 x 42', '  4    2', '  ----', '1|    |', '0|    |')\n    ,@('33 x 46', '  4    6', '  ----', '3|    |', '3|    |')\n    ,@('42 x 71', '  7    1', '  ----', '4|    |', '2|    |')\n    ,@('87 x 91', '  9    1', '  ----', '8|    |', '7|    |')\n    ,@('20 x 80', '  8    0', '  ----', '2|    |', '0|    |')\n    ,@('73 x 17', '  1    7', '  ----', '7|    |', '3|    |')\n    ,@('33 x 15', '  1    5', '  ----', '3|    |', '3|    |')\n    ,@('23 x 64', '  6    4', '  ----', '2|    |', '3|    |')\n    ,@('65 x 59', '  5    9', '  ----', '6|    |', '5|    |')\n    ,@('64 x 32', '  3    2', '  ----', '6|    |', '4|    |')\n    ,@('99 x 48', '  4    8', '  ----', '9|    |', '9|    |')\n    ,@('32 x 30', '  3    0', '  ----', '3|    |', '2|    |')\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$lineBreak = [char]11   # vertical tab == Word's manual line break (<w:br/>)\n$cellIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $lines = $afterGrid[$cellIndex]\n        $cellIndex++\n        if ($lines -eq $null) { continue }\n\n        $cell = $t.Cell($r, $c)\n        $newText = [string]::Join($lineBreak, $lines)\n        $cell.Range.Text = $newText\n    }\n}\n"}
